# Weekly update for "Hortaliza, Terminal La Palmera de La Serena - Brócoli":
# insert a new week's worth of data (2 rows: Primera / Segunda) at the top
# of the data block (row 909), pushing the rest of the table down by two
# rows and extending the used range to A1:R950.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 909 - this shifts every
# row from 909 down onward by two rows (909->911 ... 948->950), matching the
# diff exactly, and also bumps the sheet dimension to A1:R950 automatically.
$ws.Rows("909:910").Insert()

# New row 909: Brócoli, Primera, fecha 2023-01-13 (serial 44939)
$ws.Cells.Item(909, 1).Value = 8
$ws.Cells.Item(909, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(909, 3).Value = "Coquimbo"
$ws.Cells.Item(909, 4).Value = 44939
$ws.Cells.Item(909, 5).Value = 4
$ws.Cells.Item(909, 6).Value = 100112023
$ws.Cells.Item(909, 7).Value = "Brócoli"
$ws.Cells.Item(909, 8).Value = "Sin especificar"
$ws.Cells.Item(909, 9).Value = "Primera"
$ws.Cells.Item(909, 10).Value = 2600
$ws.Cells.Item(909, 11).Value = 900
$ws.Cells.Item(909, 12).Value = 1000
$ws.Cells.Item(909, 13).Value = 950
$ws.Cells.Item(909, 14).Value = "$/unidad"
$ws.Cells.Item(909, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(909, 16).Value = 950
$ws.Cells.Item(909, 17).Value = 1
$ws.Cells.Item(909, 18).Value = "Hortaliza"

# New row 910: Brócoli, Segunda, fecha 2023-01-13 (serial 44939)
$ws.Cells.Item(910, 1).Value = 8
$ws.Cells.Item(910, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(910, 3).Value = "Coquimbo"
$ws.Cells.Item(910, 4).Value = 44939
$ws.Cells.Item(910, 5).Value = 4
$ws.Cells.Item(910, 6).Value = 100112023
$ws.Cells.Item(910, 7).Value = "Brócoli"
$ws.Cells.Item(910, 8).Value = "Sin especificar"
$ws.Cells.Item(910, 9).Value = "Segunda"
$ws.Cells.Item(910, 10).Value = 1400
$ws.Cells.Item(910, 11).Value = 700
$ws.Cells.Item(910, 12).Value = 800
$ws.Cells.Item(910, 13).Value = 750
$ws.Cells.Item(910, 14).Value = "$/unidad"
$ws.Cells.Item(910, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(910, 16).Value = 750
$ws.Cells.Item(910, 17).Value = 1
$ws.Cells.Item(910, 18).Value = "Hortaliza"
